# daily auto push: 2026-01-19 06:56 UTC
#
# The daily scraper appended a new observation for 2026/01/19 13:00 that
# belongs right after the existing 2026/01/19 rows (row 680) and before the
# 2026/12/29 rows that already follow it. So: insert one row at row 681,
# shifting the old rows 681..722 down to 682..723 (and extending the used
# range from D722 to D723), then fill the new row with the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record - shifts rows 681..722 down to 682..723.
$ws.Rows.Item(681).Insert()

# Seed the new row from its neighbour (row 680, also 2026/01/19) so it
# inherits the same "date/weekday as text" cell formatting instead of
# Excel's automatic date-literal parsing, then overwrite the two columns
# that actually differ for this new reading.
$ws.Range("A680:D680").Copy()
$ws.Range("A681:D681").PasteSpecial()

$ws.Range("C681").Value = 13
$ws.Range("D681").Value = 144
